$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 269, shifting existing rows 269:367 down to 270:368
$ws.Rows(269).Insert()

# Populate the newly inserted row 269 with the new record
$ws.Range("A269").Value = 4
$ws.Range("B269").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C269").Value = "Los Lagos"
$ws.Range("D269").Value = 44900
$ws.Range("E269").Value = 10
$ws.Range("F269").Value = 100112040
$ws.Range("G269").Value = "Cilantro"
$ws.Range("H269").Value = "Sin especificar"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 70
$ws.Range("K269").Value = 6000
$ws.Range("L269").Value = 6000
$ws.Range("M269").Value = 6000
$ws.Range("N269").Value = "$/docena de atados (2 kilos)"
$ws.Range("O269").Value = "Región de La Araucanía"
$ws.Range("P269").Value = 3000
$ws.Range("Q269").Value = 2
$ws.Range("R269").Value = "Hortaliza"
